# GlyphPosition.xlsx - "Controller improved for Iphone (still not perfect)"
#
# The scale factor in Sheet2!C1 changes from 0.3 to 0.2, which ripples through
# the dependent ROUND() formulas in columns D/E and the concatenated
# "{x: .., y: ..}," strings in column G. The sheet's active selection also
# moves from G1:G11 to a single cell, I23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update the scale factor - every dependent formula (D, E, G columns)
# recalculates automatically.
$ws.Range("C1").Value = 0.2

# Leave the final selection on I23, matching the saved sheet view.
$ws.Activate()
$ws.Range("I23").Select()
